# Edit script: "Fruta / hortaliza, semanal" weekly update.
# Inserts 6 new weekly price rows for Tomate (Larga vida) at rows 1208-1213,
# shifting the previously-recorded rows 1208-1284 down to 1214-1290.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 blank rows before the current row 1208 (pushes old 1208:1284 -> 1214:1290).
$ws.Rows("1208:1213").Insert()

# Common/constant columns for every data row in this sheet.
$marketId   = 12
$marketName = "Mapocho Venta Directa de Santiago"
$region     = "Metropolitana"
$codreg     = 13
$catId      = 100112020
$categoria  = "Tomate"
$unidad     = "$/bandeja 18 kilos"
$kgUnidad   = 18
$clasif     = "Hortaliza"

# New weekly records (newest week) for rows 1208-1213.
$newRows = @(
    @{ Row=1208; Fecha=44610; Variedad="Larga vida"; Calidad="Extra";   Volumen=1040; PMin=13000; PMax=14000; PProm=13538; Origen="Provincia de Quillota";  PKg=752 }
    @{ Row=1209; Fecha=44610; Variedad="Larga vida"; Calidad="Extra";   Volumen=800;  PMin=12000; PMax=13000; PProm=12450; Origen="Región de O'Higgins";    PKg=692 }
    @{ Row=1210; Fecha=44610; Variedad="Larga vida"; Calidad="Primera"; Volumen=1080; PMin=10000; PMax=11000; PProm=10444; Origen="Provincia de Quillota";  PKg=580 }
    @{ Row=1211; Fecha=44610; Variedad="Larga vida"; Calidad="Primera"; Volumen=980;  PMin=9500;  PMax=10000; PProm=9735;  Origen="Región de O'Higgins";    PKg=541 }
    @{ Row=1212; Fecha=44610; Variedad="Larga vida"; Calidad="Segunda"; Volumen=860;  PMin=7000;  PMax=8000;  PProm=7535;  Origen="Provincia de Quillota";  PKg=419 }
    @{ Row=1213; Fecha=44610; Variedad="Larga vida"; Calidad="Segunda"; Volumen=540;  PMin=6500;  PMax=7000;  PProm=6741;  Origen="Región de O'Higgins";    PKg=374 }
)

foreach ($rec in $newRows) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value  = $marketId
    $ws.Cells.Item($r, 2).Value  = $marketName
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $rec.Fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $catId
    $ws.Cells.Item($r, 7).Value  = $categoria
    $ws.Cells.Item($r, 8).Value  = $rec.Variedad
    $ws.Cells.Item($r, 9).Value  = $rec.Calidad
    $ws.Cells.Item($r, 10).Value = $rec.Volumen
    $ws.Cells.Item($r, 11).Value = $rec.PMin
    $ws.Cells.Item($r, 12).Value = $rec.PMax
    $ws.Cells.Item($r, 13).Value = $rec.PProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $rec.Origen
    $ws.Cells.Item($r, 16).Value = $rec.PKg
    $ws.Cells.Item($r, 17).Value = $kgUnidad
    $ws.Cells.Item($r, 18).Value = $clasif
}
